$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# New header row (row 1) -- columns F..T
# ---------------------------------------------------------------
$ws.Cells.Item(1, 6).Value  = "name"
$ws.Cells.Item(1, 7).Value  = "categ_id"
$ws.Cells.Item(1, 8).Value  = "type"
$ws.Cells.Item(1, 9).Value  = "*TMPL*sale_ok"
$ws.Cells.Item(1, 10).Value = "*TMPL*purchase_ok"
$ws.Cells.Item(1, 11).Value = "standard_price"
$ws.Cells.Item(1, 12).Value = "*TMPL*list_price"
$ws.Cells.Item(1, 13).Value = "*TMPL*invoice_policy"
$ws.Cells.Item(1, 14).Value = "*TMPL*image"
$ws.Cells.Item(1, 15).Value = "image_variant"
$ws.Cells.Item(1, 16).Value = "*TMPL*description_sale"
$ws.Cells.Item(1, 17).Value = "*TMPL*description_purchase"
$ws.Cells.Item(1, 18).Value = "uom_id"
$ws.Cells.Item(1, 19).Value = "uom_po_id"
$ws.Cells.Item(1, 20).Value = "*TMPL*product_brand_id"

# categ_id header column uses the text ("@") number format, same as
# the existing "attribute: Size" style already present in the sheet.
$ws.Cells.Item(1, 7).NumberFormat = "@"

# ---------------------------------------------------------------
# New data row (row 2) -- columns F..T
# ---------------------------------------------------------------
$ws.Cells.Item(2, 6).Value  = "Template 1"
$ws.Cells.Item(2, 7).Value  = "Categ 1"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value  = "service"
# NOTE: assigning the literal strings "True"/"False" via .Value makes the
# engine auto-coerce them into real Booleans (t="b"), same as Excel does
# on direct input. The source file needs them stored as plain text
# (t="s"), so they are entered with a leading apostrophe (forces text)
# and then the "quote prefix" styling that apostrophe-entry adds is
# stripped by pasting the (unstyled) format of A2 back onto the cells.
$ws.Cells.Item(2, 9).Value  = "'True"
$ws.Cells.Item(2, 10).Value = "'False"
$ws.Cells.Item(2, 1).Copy()
$ws.Range("I2:J2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(2, 11).Value = 10.5
$ws.Cells.Item(2, 12).Value = 15.99
$ws.Cells.Item(2, 13).Value = "delivery"
$ws.Cells.Item(2, 14).Value = "https://trey.es/web/image/website/1/logo?unique=20ee16c"
$ws.Cells.Item(2, 16).Value = "Description for customers."
$ws.Cells.Item(2, 17).Value = "Description for suppliers."
$ws.Cells.Item(2, 18).Value = "Unit(s)"
$ws.Cells.Item(2, 19).Value = "Unit(s)"
$ws.Cells.Item(2, 20).Value = "Brand 1"

# ---------------------------------------------------------------
# Column width: column N (14) gets a lot wider to fit the image URL
# ---------------------------------------------------------------
$ws.Columns.Item(14).ColumnWidth = 68.2166666666667

# ---------------------------------------------------------------
# Sheet view: scroll so column G is the leftmost visible column and
# select N28 as the active cell (matches the authored selection).
# ---------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 7
[void]$ws.Range("N28").Select()
